# Add 2024-season NFL game rows to the "games" sheet, matching the upload
# that extended the table from row 545 down to row 773, and re-point the
# active sheet/selection to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("games")

# ---------------------------------------------------------------------
# 1) Finish the two trailing 2023-season rows that had only A:C filled
#    in (the H column already carried a shared "=E<row>" formula whose
#    value was showing 0 because D:G were still blank).
# ---------------------------------------------------------------------
$ws.Range("D544").Value = "DAL"
$ws.Range("E544").Value = "WAS"
$ws.Range("F544").Value = 38
$ws.Range("G544").Value = 10
$ws.Range("H544").Formula = "=E544"

$ws.Range("D545").Value = "BUF"
$ws.Range("E545").Value = "MIA"
$ws.Range("F545").Value = 21
$ws.Range("G545").Value = 14
$ws.Range("H545").Formula = "=E545"

# ---------------------------------------------------------------------
# 2) Append the 2024 week-1 games (rows 546-558), fully populated.
#    Grab the date format (m/d/yyyy) from an existing date cell first so
#    the new C column cells reuse the same style record instead of a new
#    one.
# ---------------------------------------------------------------------
$ws.Range("C543").Copy()
$ws.Range("C546:C561").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$games2024 = @(
  @(546, 1, 45540, "BAL", "KC",  20, 27),
  @(547, 1, 45541, "GB",  "PHI", 29, 34),
  @(548, 1, 45543, "CAR", "NO",  10, 47),
  @(549, 1, 45543, "TEN", "CHI", 17, 24),
  @(550, 1, 45543, "NE",  "CIN", 16, 10),
  @(551, 1, 45543, "PIT", "ATL", 18, 10),
  @(552, 1, 45543, "ARI", "BUF", 28, 34),
  @(553, 1, 45543, "MIN", "NYG", 28, 6),
  @(554, 1, 45543, "JAX", "MIA", 17, 20),
  @(555, 1, 45543, "HOU", "IND", 29, 27),
  @(556, 1, 45543, "LV",  "LAC", 10, 22),
  @(557, 1, 45543, "DEN", "SEA", 20, 26),
  @(558, 1, 45543, "WAS", "TB",  20, 37)
)

foreach ($g in $games2024) {
    $r = $g[0]
    $ws.Cells.Item($r, 1).Value = 2024
    $ws.Cells.Item($r, 2).Value = $g[1]
    $ws.Cells.Item($r, 3).Value = $g[2]
    $ws.Cells.Item($r, 4).Value = $g[3]
    $ws.Cells.Item($r, 5).Value = $g[4]
    $ws.Cells.Item($r, 6).Value = $g[5]
    $ws.Cells.Item($r, 7).Value = $g[6]
    $ws.Cells.Item($r, 8).Formula = "=E" + $r
}

# ---------------------------------------------------------------------
# 3) A couple more rows that only got season/week/date filled in before
#    the paste of the rest of the 2024 schedule trailed off (559-561),
#    followed by a long run of rows where only the season (2024) made
#    it into column A (562-773).
# ---------------------------------------------------------------------
$partialDateRows = @(
  @(559, 1, 45543),
  @(560, 1, 45543),
  @(561, 1, 45544)
)
foreach ($p in $partialDateRows) {
    $r = $p[0]
    $ws.Cells.Item($r, 1).Value = 2024
    $ws.Cells.Item($r, 2).Value = $p[1]
    $ws.Cells.Item($r, 3).Value = $p[2]
}

for ($r = 562; $r -le 773; $r++) {
    $ws.Cells.Item($r, 1).Value = 2024
}

# ---------------------------------------------------------------------
# 4) Leave the view the way the author left it: "games" active/selected
#    with D559 highlighted and scrolled near the bottom of the new data;
#    "2025 schedule" no longer the selected tab.
# ---------------------------------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("D559").Select()
